# "Generate Report for Handback" - localization-status.xlsx update
#
# Summary of the change being applied:
#  - Overview!Status text "Ready for handoff" -> "Handed back: in sync with en-US"
#    (shared by Overview E2/F2/E3/F3 and the per-locale sheets' Status column C2/C3)
#  - zh-cn / de-de sheets: "Latest Target File" (col I) and "Latest Handback File"
#    (col J) rows 2/3 are now populated (were blank) because a handback report was
#    generated. Col I gets a hyperlink to the source .md file (like col A) and col J
#    gets the generated handback xliff file name.
#  - de-de sheet: "Latest Handback DateTime" (col K) rows 2/3 move from the zero-date
#    placeholder to the real handback timestamp.
#  - Column widths for the newly-populated / widened columns grow to fit content.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# ----------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Overview column widths (Status columns widen to fit the longer text)
$overview.Range("E1").ColumnWidth = 30
$overview.Range("F1").ColumnWidth = 30

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

# Hyperlink-style text color (matches the existing custom "HyperLink" cell style:
# font color FF6495ED / RGB 100,149,237, single underline)
$hyperlinkColor = 15570276  # OLE BGR encoding of RGB(100,149,237)

# ----------------------------------------------------------------------------
# 2) zh-cn sheet: populate "Latest Target File" (I) / "Latest Handback File" (J)
# ----------------------------------------------------------------------------
$zhcn.Hyperlinks.Add(
    $zhcn.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54d5c65db63870f2fceeed6f0bcebdc693067ce6/e2e/1e851055-11df-4ef1-924f-439760840548.md",
    "",
    "",
    "1e851055-11df-4ef1-924f-439760840548.md")
$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = $hyperlinkColor
$zhcn.Range("J2").Value = "1e851055-11df-4ef1-924f-439760840548.d72fbf889ad9b831bb496ae3d0e73ee957d86ff7.zh-cn.xlf"

$zhcn.Hyperlinks.Add(
    $zhcn.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54d5c65db63870f2fceeed6f0bcebdc693067ce6/e2e/8c02821f-2870-4848-b5d5-3f80b9d200ad.md",
    "",
    "",
    "8c02821f-2870-4848-b5d5-3f80b9d200ad.md")
$zhcn.Range("I3").Font.Underline = 2
$zhcn.Range("I3").Font.Color = $hyperlinkColor
$zhcn.Range("J3").Value = "8c02821f-2870-4848-b5d5-3f80b9d200ad.2486bf1c9fc10e2f13e8e658ad0dc16394189bc6.zh-cn.xlf"

# zh-cn column widths: Status (C) widens, Latest Target File (I) / Latest Handback
# File (J) widen to fit the newly-populated long file names
$zhcn.Range("C1").ColumnWidth = 30
$zhcn.Range("I1").ColumnWidth = 40
$zhcn.Range("J1").ColumnWidth = 40

# ----------------------------------------------------------------------------
# 3) de-de sheet: populate "Latest Target File" (I) / "Latest Handback File" (J)
#    and the "Latest Handback DateTime" (K)
# ----------------------------------------------------------------------------
$dede.Hyperlinks.Add(
    $dede.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54d5c65db63870f2fceeed6f0bcebdc693067ce6/e2e/1e851055-11df-4ef1-924f-439760840548.md",
    "",
    "",
    "1e851055-11df-4ef1-924f-439760840548.md")
$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = $hyperlinkColor
$dede.Range("J2").Value = "1e851055-11df-4ef1-924f-439760840548.d72fbf889ad9b831bb496ae3d0e73ee957d86ff7.de-de.xlf"
$dede.Range("K2").Value = "2016-09-04 15:05:21"

$dede.Hyperlinks.Add(
    $dede.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54d5c65db63870f2fceeed6f0bcebdc693067ce6/e2e/8c02821f-2870-4848-b5d5-3f80b9d200ad.md",
    "",
    "",
    "8c02821f-2870-4848-b5d5-3f80b9d200ad.md")
$dede.Range("I3").Font.Underline = 2
$dede.Range("I3").Font.Color = $hyperlinkColor
$dede.Range("J3").Value = "8c02821f-2870-4848-b5d5-3f80b9d200ad.2486bf1c9fc10e2f13e8e658ad0dc16394189bc6.de-de.xlf"
$dede.Range("K3").Value = "2016-09-04 15:05:21"

# de-de column widths: same widening pattern as zh-cn
$dede.Range("C1").ColumnWidth = 30
$dede.Range("I1").ColumnWidth = 40
$dede.Range("J1").ColumnWidth = 40

Write-Host "Handback report generated."
